# master-machine_type.xlsx — "Adding Master Data XLS"
#
# The sheet held sample master-data rows for machine type "DNG / Dongle".
# The commit swaps that sample row set for a "DKS / Dekstop" (Desktop
# Computer) set, keeping the same 7-column schema
# (code, name, descr, lang_code, is_active, cr_by, cr_dtimes) and the same
# 3 language rows (eng, ara, fra) per machine type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (eng): DNG/Dongle/"To run enrollment client" -> DKS/Dekstop/"Desktop Computer"
$ws.Range("A2").Value2 = "DKS"
$ws.Range("B2").Value2 = "Dekstop"
$ws.Range("C2").Value2 = "Desktop Computer"
$ws.Range("D2").Value2 = "eng"

# Row 3 (ara): DNG/دونجل/لتشغيل عميل التسجيل -> DKS/الحاسوب/أجهزة الكمبيوتر المكتبية
$ws.Range("A3").Value2 = "DKS"
$ws.Range("B3").Value2 = "الحاسوب"
$ws.Range("C3").Value2 = "أجهزة الكمبيوتر المكتبية"
$ws.Range("D3").Value2 = "ara"

# Row 4 (fra): DNG/Dongle/"Pour exécuter le client dinscription" -> DKS/Ordinateur/"Ordinateurs de bureau"
$ws.Range("A4").Value2 = "DKS"
$ws.Range("B4").Value2 = "Ordinateur"
$ws.Range("C4").Value2 = "Ordinateurs de bureau"
$ws.Range("D4").Value2 = "fra"

# is_active / cr_by / cr_dtimes (E:G) are unchanged (TRUE, superadmin, now())
# across all three rows — left untouched.

# The saved file's selection moved to D10 before the next save.
$ws.Range("D10").Select() | Out-Null

# Page setup was touched (paper size A4, portrait) which is what produced
# the new <pageSetup .../> element on the sheet.
$ps = $ws.PageSetup
$ps.PaperSize = 9        # xlPaperA4
$ps.Orientation = 1      # xlPortrait
